$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme
Write-Host "Count: $($tcs.Count)"
for ($i=1; $i -le $tcs.Count; $i++) {
    try {
        $c = $tcs.Item($i)
        Write-Host "Item $i : $c  RGB=$($c.RGB)"
    } catch {
        Write-Host "Item $i ERR: $_"
    }
}
